$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = -0.4000000059604645
$ws.Range("R3").Value = -1.600000023841858
$ws.Range("X3").Value = -30.5

$ws.Range("E4").Value = -1.5
$ws.Range("K4").Value = -0.699999988079071
$ws.Range("R4").Value = -1.600000023841858
$ws.Range("X4").Value = -18.60000038146973

$ws.Range("E5").Value = 0
$ws.Range("K5").Value = -0.300000011920929
$ws.Range("R5").Value = -2
$ws.Range("X5").Value = -25.10000038146973

$ws.Range("E6").Value = 0
$ws.Range("K6").Value = -3.900000095367432
$ws.Range("R6").Value = -2
$ws.Range("X6").Value = -24.10000038146973

$ws.Range("K7").Value = -1.700000047683716
$ws.Range("R7").Value = -4.300000190734863
$ws.Range("X7").Value = -26.29999923706055

$ws.Range("E8").Value = -0.300000011920929
$ws.Range("K8").Value = -0.4000000059604645
$ws.Range("X8").Value = -23.39999961853027

$ws.Range("E9").Value = -2.400000095367432
$ws.Range("K9").Value = -0.1000000014901161
$ws.Range("R9").Value = -2.5
$ws.Range("X9").Value = -30.70000076293945

$ws.Range("E10").Value = -2.599999904632568
$ws.Range("R10").Value = -2.5
$ws.Range("X10").Value = -26.39999961853027

$ws.Range("E11").Value = -0.1000000014901161
$ws.Range("K11").Value = -1.299999952316284
$ws.Range("R11").Value = -3.099999904632568
$ws.Range("X11").Value = -26.60000038146973

$ws.Range("E12").Value = -0.699999988079071
$ws.Range("R12").Value = -3.099999904632568
$ws.Range("X12").Value = -30

$ws.Range("X13").Value = -59.90000152587891

$ws.Range("E14").Value = -3.299999952316284
$ws.Range("K14").Value = -0.2000000029802322
$ws.Range("R14").Value = -4
$ws.Range("X14").Value = -29

$ws.Range("R15").Value = -4.699999809265137
$ws.Range("X15").Value = -39.70000076293945
